# Update "想去人数" (want-to-go count) figures and one refreshed cover-image
# URL across the four sheets, matching the upstream data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibitions) -------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value  = 2683
$ws1.Range("F3").Value  = 1034
$ws1.Range("F4").Value  = 19249
$ws1.Range("F6").Value  = 2159
$ws1.Range("F7").Value  = 728
$ws1.Range("F9").Value  = 410
$ws1.Range("F10").Value = 668
$ws1.Range("F11").Value = 224
$ws1.Range("F12").Value = 238
$ws1.Range("F14").Value = 349
$ws1.Range("F16").Value = 247
$ws1.Range("F18").Value = 170
$ws1.Range("F19").Value = 14
$ws1.Range("F21").Value = 88

# --- Sheet 2: 演出 (Performances) -------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value  = 184
$ws2.Range("I5").Value  = "//i0.hdslb.com/bfs/openplatform/202409/lyCquRtq1727079420725.png"
$ws2.Range("F7").Value  = 265
$ws2.Range("F15").Value = 57

# --- Sheet 3: 本地生活 (Local Life) -----------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 5967
$ws3.Range("F3").Value = 623

# --- Sheet 4: 全部类型 (All Types) ------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value  = 5967
$ws4.Range("F3").Value  = 623
$ws4.Range("F5").Value  = 184
$ws4.Range("F7").Value  = 2683
$ws4.Range("F8").Value  = 1034
$ws4.Range("F9").Value  = 19249
$ws4.Range("I11").Value = "//i0.hdslb.com/bfs/openplatform/202409/lyCquRtq1727079420725.png"
$ws4.Range("F14").Value = 265
$ws4.Range("F15").Value = 2159
$ws4.Range("F16").Value = 728
$ws4.Range("F19").Value = 410
$ws4.Range("F20").Value = 668
$ws4.Range("F21").Value = 224
$ws4.Range("F22").Value = 238
$ws4.Range("F27").Value = 349
$ws4.Range("F30").Value = 247
$ws4.Range("F34").Value = 170
$ws4.Range("F35").Value = 57
$ws4.Range("F37").Value = 14
$ws4.Range("F45").Value = 88
